# Deploy the implementation guide.
#
# Two changes:
#  1. "Metadata" sheet: bump the generation Date from
#     2023-01-30T21:30:05+00:00 to 2023-02-09T16:19:35+00:00.
#  2. "Concepts" sheet: append a new concept row (row 10) for code "RGDI+"
#     - column A (Level) mirrors the "1" used by every other row
#     - columns B (Code) and C (Display) both hold "RGDI+"
#     - column D (Definition) stays empty, same as the existing rows
#     Copying the last data row down preserves its style/formatting for
#     the new row, matching the rest of the table.

$wb = $excel.ActiveWorkbook

# 1. Update the Date value on the Metadata sheet.
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2023-02-09T16:19:35+00:00"

# 2. Add the new "RGDI+" concept row on the Concepts sheet.
$concepts = $wb.Worksheets.Item("Concepts")
$concepts.Range("A9:D9").Copy($concepts.Range("A10:D10"))
$concepts.Range("B10").Value = "RGDI+"
$concepts.Range("C10").Value = "RGDI+"
